# Dataset.xlsx update — 12/18/2017 changes
# - URL sheet: move selection from F3 to D3
# - Disruptions sheet: move selection from A7 to D2, and refresh the sample
#   flight-number data (5x4 grid of FlightSort columns) with the new
#   2017-12-15 flight ids, matching the 15/60/180-min Summary Drawer +
#   Filter view-options dataset refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# URL sheet — selection only
# ---------------------------------------------------------------------
$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------
# Disruptions sheet — new flight data + selection
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Disruptions")
$ws.Activate() | Out-Null

# Columns C2:C5 and E2 currently carry style index 6 / 4 (General number
# format); the refreshed sheet re-enters them with the Text-formatted
# style already used by the B/D columns (index 32). Copy that format
# across first so the re-typed values land in the same style.
$donor = $ws.Range("B2")
$donor.Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

# Write the new flight numbers — new values first (in the order the new
# 2017-12-15 flight ids first appear), then the remaining cells that
# reuse either a new value or one of the still-present FlightSort*
# labels, so the shared-string table comes out in the same order/shape
# as the refreshed workbook.
$ws.Range("B2").Value = "FL-ZZ-158-20171215-DCA-LAX-0"
$ws.Range("B3").Value = "FL-ZZ-336-20171215-GSP-LAX-0"
$ws.Range("B4").Value = "FL-ZZ-876-20171215-PNM-LAX-0"
$ws.Range("B5").Value = "FL-ZZ-814-20171215-CLT-LAX-0"
$ws.Range("E2").Value = "FL-ZZ-60-20171215-PIE-LAX-0"
$ws.Range("C2").Value = "FL-ZZ-144-20171215-SEA-LAX-0"
$ws.Range("C4").Value = "FL-ZZ-26-20171215-PHL-LAX-0"
$ws.Range("D2").Value = "FL-ZZ-886-20171215-SDF-LAX-0"
$ws.Range("D4").Value = "FL-ZZ-868-20171215-PIT-LAX-0"
$ws.Range("E3").Value = "FL-ZZ-64-20171215-PIE-LAX-0"
$ws.Range("C5").Value = "FL-ZZ-28-20171215-PHL-LAX-0"

$ws.Range("C1").Value = "FlightSortAsc"
$ws.Range("D1").Value = "FlightSortDesc"
$ws.Range("E1").Value = "FlightSortNewAsc"
$ws.Range("C3").Value = "FL-ZZ-158-20171215-DCA-LAX-0"
$ws.Range("D3").Value = "FL-ZZ-876-20171215-PNM-LAX-0"
$ws.Range("E4").Value = "FL-ZZ-144-20171215-SEA-LAX-0"
$ws.Range("D5").Value = "FL-ZZ-814-20171215-CLT-LAX-0"
$ws.Range("E5").Value = "FL-ZZ-886-20171215-SDF-LAX-0"
$ws.Range("B6").Value = "FL-ZZ-60-20171215-PIE-LAX-0"
$ws.Range("B7").Value = "FL-ZZ-26-20171215-PHL-LAX-0"

$ws.Range("D2").Select() | Out-Null
